$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("A2").Value = "ECs"
$ws.Range("B2").Value = "Efnb1"
$ws.Range("C2").Value = "Ephb1"
$ws.Range("D2").Value = "ECs"
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 7.079689000000001
$ws.Range("H2").Value = 21.239067
$ws.Range("I2").Value = 0.5033576067109902
$ws.Range("J2").Value = 0.5033576067109902
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 1.916802666666667
$ws.Range("N2").Value = 5.750408
$ws.Range("O2").Value = 0.3865473586068074
$ws.Range("P2").Value = 0.3865473586068074
$ws.Range("Q2").Value = 13.57036675437067
$ws.Range("R2").Value = 122.133300789336
$ws.Range("S2").Value = 0.1945715533087774
$ws.Range("T2").Value = 0.1945715533087775

# Row 3
$ws.Range("A3").Value = "ECs"
$ws.Range("B3").Value = "Efnb1"
$ws.Range("C3").Value = "Ephb1"
$ws.Range("D3").Value = "sCs"
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 7.079689000000001
$ws.Range("H3").Value = 21.239067
$ws.Range("I3").Value = 0.5033576067109902
$ws.Range("J3").Value = 0.5033576067109902
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 3.041975666666666
$ws.Range("N3").Value = 9.125926999999999
$ws.Range("O3").Value = 0.6134526413931926
$ws.Range("P3").Value = 0.6134526413931926
$ws.Range("Q3").Value = 21.53624166556767
$ws.Range("R3").Value = 193.826174990109
$ws.Range("S3").Value = 0.3087860534022127
$ws.Range("T3").Value = 0.3087860534022127

# Row 4
$ws.Range("A4").Value = "FAPs"
$ws.Range("B4").Value = "Efnb1"
$ws.Range("C4").Value = "Ephb1"
$ws.Range("D4").Value = "ECs"
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 4.058683666666667
$ws.Range("H4").Value = 12.176051
$ws.Range("I4").Value = 0.2885676612136944
$ws.Range("J4").Value = 0.2885676612136945
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 1.916802666666667
$ws.Range("N4").Value = 5.750408
$ws.Range("O4").Value = 0.3865473586068074
$ws.Range("P4").Value = 0.3865473586068074
$ws.Range("Q4").Value = 7.779695675423111
$ws.Range("R4").Value = 70.017261078808
$ws.Range("S4").Value = 0.1115450672214976
$ws.Range("T4").Value = 0.1115450672214977

# Row 5
$ws.Range("A5").Value = "FAPs"
$ws.Range("B5").Value = "Efnb1"
$ws.Range("C5").Value = "Ephb1"
$ws.Range("D5").Value = "sCs"
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 4.058683666666667
$ws.Range("H5").Value = 12.176051
$ws.Range("I5").Value = 0.2885676612136944
$ws.Range("J5").Value = 0.2885676612136945
$ws.Range("K5").Value = 3
$ws.Range("L5").Value = 1
$ws.Range("M5").Value = 3.041975666666666
$ws.Range("N5").Value = 9.125926999999999
$ws.Range("O5").Value = 0.6134526413931926
$ws.Range("P5").Value = 0.6134526413931926
$ws.Range("Q5").Value = 12.34641695269744
$ws.Range("R5").Value = 111.117752574277
$ws.Range("S5").Value = 0.1770225939921968
$ws.Range("T5").Value = 0.1770225939921968

# Row 6
$ws.Range("A6").Value = "M2"
$ws.Range("B6").Value = "Efnb1"
$ws.Range("C6").Value = "Ephb1"
$ws.Range("D6").Value = "ECs"
$ws.Range("E6").Value = 3
$ws.Range("F6").Value = 1
$ws.Range("G6").Value = 0.463315
$ws.Range("H6").Value = 1.389945
$ws.Range("I6").Value = 0.03294115455541936
$ws.Range("J6").Value = 0.03294115455541936
$ws.Range("K6").Value = 2
$ws.Range("L6").Value = 0.6666666666666666
$ws.Range("M6").Value = 1.916802666666667
$ws.Range("N6").Value = 5.750408
$ws.Range("O6").Value = 0.3865473586068074
$ws.Range("P6").Value = 0.3865473586068074
$ws.Range("Q6").Value = 0.8880834275066666
$ws.Range("R6").Value = 7.99275084756
$ws.Range("S6").Value = 0.01273331628285596
$ws.Range("T6").Value = 0.01273331628285596

# Row 7
$ws.Range("A7").Value = "M2"
$ws.Range("B7").Value = "Efnb1"
$ws.Range("C7").Value = "Ephb1"
$ws.Range("D7").Value = "sCs"
$ws.Range("E7").Value = 3
$ws.Range("F7").Value = 1
$ws.Range("G7").Value = 0.463315
$ws.Range("H7").Value = 1.389945
$ws.Range("I7").Value = 0.03294115455541936
$ws.Range("J7").Value = 0.03294115455541936
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 3.041975666666666
$ws.Range("N7").Value = 9.125926999999999
$ws.Range("O7").Value = 0.6134526413931926
$ws.Range("P7").Value = 0.6134526413931926
$ws.Range("Q7").Value = 1.409392956001666
$ws.Range("R7").Value = 12.684536604015
$ws.Range("S7").Value = 0.0202078382725634
$ws.Range("T7").Value = 0.0202078382725634

# Row 8
$ws.Range("A8").Value = "sCs"
$ws.Range("B8").Value = "Efnb1"
$ws.Range("C8").Value = "Ephb1"
$ws.Range("D8").Value = "ECs"
$ws.Range("E8").Value = 3
$ws.Range("F8").Value = 1
$ws.Range("G8").Value = 2.463241333333333
$ws.Range("H8").Value = 7.389724
$ws.Range("I8").Value = 0.175133577519896
$ws.Range("J8").Value = 0.175133577519896
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 1.916802666666667
$ws.Range("N8").Value = 5.750408
$ws.Range("O8").Value = 0.3865473586068074
$ws.Range("P8").Value = 0.3865473586068074
$ws.Range("Q8").Value = 4.721547556376889
$ws.Range("R8").Value = 42.493928007392
$ws.Range("S8").Value = 0.06769742179367633
$ws.Range("T8").Value = 0.06769742179367634

# Row 9
$ws.Range("A9").Value = "sCs"
$ws.Range("B9").Value = "Efnb1"
$ws.Range("C9").Value = "Ephb1"
$ws.Range("D9").Value = "sCs"
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = 1
$ws.Range("G9").Value = 2.463241333333333
$ws.Range("H9").Value = 7.389724
$ws.Range("I9").Value = 0.175133577519896
$ws.Range("J9").Value = 0.175133577519896
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 1
$ws.Range("M9").Value = 3.041975666666666
$ws.Range("N9").Value = 9.125926999999999
$ws.Range("O9").Value = 0.6134526413931926
$ws.Range("P9").Value = 0.6134526413931926
$ws.Range("Q9").Value = 7.493120197127555
$ws.Range("R9").Value = 67.43808177414799
$ws.Range("S9").Value = 0.1074361557262196
$ws.Range("T9").Value = 0.1074361557262196
